# Update "Ready for handoff" status cells to "In Translation" across all
# worksheets that reference the shared string, and narrow the
# "Latest HO Xliff Generate Date"/"Status" columns that were auto-fit to a
# shorter width ("17.22" -> "13.41" character units) when the status text
# changed from "Ready for handoff" to the shorter "In Translation".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Overview ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Columns E (zh-cn) and F (de-de) shrink to match the new, shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.57
$wsOverview.Columns.Item(6).ColumnWidth = 12.57

# --- Sheet 2: zh-cn --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.57

# --- Sheet 3: de-de --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.57
